$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 161 (pushes existing rows 161-249 down to 162-250).
$ws.Rows(161).Insert()

# Populate the newly inserted row 161 with the new record.
$ws.Range("A161").Value2 = 10
$ws.Range("B161").Value2 = "Vega Modelo de Temuco"
$ws.Range("C161").Value2 = "La Araucanía"
$ws.Range("D161").Value2 = 45097
$ws.Range("E161").Value2 = 9
$ws.Range("F161").Value2 = "Fruta"
$ws.Range("G161").Value2 = 100104
$ws.Range("H161").Value2 = "Frutos de pepita"
$ws.Range("I161").Value2 = 100104001
$ws.Range("J161").Value2 = "Granada"
$ws.Range("K161").Value2 = "Wonderfull"
$ws.Range("L161").Value2 = "Primera"
$ws.Range("M161").Value2 = 110
$ws.Range("N161").Value2 = 13000
$ws.Range("O161").Value2 = 13000
$ws.Range("P161").Value2 = 13000
$ws.Range("Q161").Value2 = "`$/bandeja 10 kilos"
$ws.Range("R161").Value2 = "Provincia de Limarí"
$ws.Range("S161").Value2 = 1300
$ws.Range("T161").Value2 = 10
